$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("January 2014")

# Update selection to C14 (matches sheetView selection in diff)
$ws.Range("C14").Select()

# Row 12: new timesheet entry
$ws.Range("A12").Value = [DateTime]"2014-01-14"
$ws.Range("B12").Value = "Resource Edit"
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = "Riaan Bekker"

# Row 13: new timesheet entry
$ws.Range("A13").Value = [DateTime]"2014-01-16"
$ws.Range("B13").Value = "Resource Upload changes and Resource Edit"
$ws.Range("C13").Value = 4.5
$ws.Range("D13").Value = "Riaan Bekker"

$wb.Save()
